$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 28
$ws.Cells.Item($row, 1).Value = "test"
$ws.Cells.Item($row, 2).Value = 1741
$ws.Cells.Item($row, 3).Value = 1851
$ws.Cells.Item($row, 4).Value = 7
$ws.Cells.Item($row, 5).Value = 1633

$textRange = $ws.Range("F28:J28")
$textRange.NumberFormat = "@"
$ws.Cells.Item($row, 6).Value = "33.28"
$ws.Cells.Item($row, 7).Value = "35.38"
$ws.Cells.Item($row, 8).Value = "0.13"
$ws.Cells.Item($row, 9).Value = "31.21"
$ws.Cells.Item($row, 10).Value = "2025-09-01 20:29:32"
$textRange.ClearFormats()
